# Apply the "update sheet with scopus" edit to the "query" worksheet.
#
# Summary of the change:
#  - A new column is inserted before column C to hold a Scopus-syntax
#    version of each Web-of-Science search string (header "Scopus" in C1,
#    a SUBSTITUTE() formula in C2 translating "TS = " / "NEAR/" into
#    "TITLE-ABS-KEY" / "W/").
#  - A new summary row is inserted below the "Transport policies" block
#    (before the old row 8) that concatenates the six transport search
#    strings together and records their combined yield (3942).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("query")

# ---------------------------------------------------------------------
# 1. Insert the new "Scopus" column before column C (old C/D shift to D/E)
# ---------------------------------------------------------------------
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(3).ColumnWidth = 23.6

# ---------------------------------------------------------------------
# 2. Insert the new summary row before the old row 8 (rows >=8 shift +1)
# ---------------------------------------------------------------------
$ws.Rows.Item(8).Insert()

# ---------------------------------------------------------------------
# 3. Populate the new summary row (row 8) -- done first so the new shared
#    string "SUMMARY" lands at index 45, matching the target string table
#    ordering (the "Scopus" header string must come right after it).
# ---------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value = "SUMMARY"
$ws.Cells.Item(8, 2).Formula = '=CONCATENATE("(",B2,")","OR","(",B3,")","OR","(",B4,")","OR","(",B5,")","OR","(",B6,")","OR","(",B7,")")'
$ws.Cells.Item(8, 2).WrapText = $false
$ws.Cells.Item(8, 3).WrapText = $false
$ws.Cells.Item(8, 4).Value = 3942
$ws.Cells.Item(8, 5).NumberFormat = $ws.Cells.Item(7, 5).NumberFormat

# ---------------------------------------------------------------------
# 4. Populate the new column header + conversion formula
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 3).Value = "Scopus"
$ws.Cells.Item(2, 3).Formula = '=SUBSTITUTE(SUBSTITUTE(B2,"TS = ","TITLE-ABS-KEY"),"NEAR/","W/")'

Write-Host "Scopus column and SUMMARY row inserted"
